# Add "compress_time" / "total_time" columns with data to the "fftss" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fftss")

$ws.Range("D14").Value = "compress_time"
$ws.Range("E14").Value = "total_time"

$ws.Range("D15").Value = 0.022270000000000002
$ws.Range("E15").Value = 0.028202000000000001

$ws.Range("D16").Value = 0.019782000000000001
$ws.Range("E16").Value = 0.022936000000000002

$ws.Range("D17").Value = 0.018460000000000001
$ws.Range("E17").Value = 0.024346

$ws.Range("D18").Value = 0.026015
$ws.Range("E18").Value = 0.030253

$ws.Range("D19").Formula = "=SUM(D15:D18)"
$ws.Range("E19").Formula = "=SUM(E15:E18)"

$ws.Range("E19").Select()
